# Update automatico via Actualizar 02-15-2021 12-39-21
#
# The worksheet keeps a rolling log of availability checks: every run
# writes a fresh timestamp into column D for the newest block of rows
# (2-15) while the older blocks (16-29 and 30-43) keep their own,
# progressively older, timestamps - they simply "age" one slot down on
# each update. This edit pushes a brand-new check timestamp into the
# newest block and shifts the two existing timestamps down into the
# next-older blocks, exactly like the previous run's history would.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest block (rows 2-15): gets the new "just checked" timestamp.
$ws.Range("D2:D15").Value = 44242.52724028876

# Middle block (rows 16-29): inherits the timestamp that used to belong
# to the newest block.
$ws.Range("D16:D29").Value = 44242.50600597222

# Oldest kept block (rows 30-43): inherits the timestamp that used to
# belong to the middle block (the previous oldest entry rolls off).
$ws.Range("D30:D43").Value = 44242.48478574074
